$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Full results" ---
$ws1 = $wb.Worksheets.Item("Full results")

# Row 2 (NULL MODEL / income)
$ws1.Range("C2").Value = 0.933934515757714
$ws1.Range("D2").Value = 0.0661353206243123
$ws1.Range("E2").Value = 1.00006983638203
$ws1.Range("J2").Value = 0.0661307022953231
$ws1.Range("K2").Value = 0.0659627311785311
$ws1.Range("L2").Value = 0.0530047004353736
$ws1.Range("M2").Value = 0.0857567329439641
$ws1.Range("N2").Value = 0.118967431613905

# Row 3 (CONDITIONAL MODEL / income)
$ws1.Range("F3").Value = 0.90118019596567
$ws1.Range("G3").Value = 0.0659673377770252

# Row 4 (COMPLETE MODEL / income)
$ws1.Range("H4").Value = 0.848171793873787
$ws1.Range("I4").Value = 0.00196805124185963
$ws1.Range("O4").Value = 0.151887435239287

# --- Sheet 2: "For plotting" ---
$ws2 = $wb.Worksheets.Item("For plotting")

# Row 2 (condind / income)
$ws2.Range("C2").Value = 0.0661307022953231
$ws2.Range("D2").Value = -0.00190399304253661
$ws2.Range("E2").Value = 0.134165397633183
$ws2.Range("F2").Value = 948

# Row 3 (completeind / income)
$ws2.Range("C3").Value = 0.118967431613905
$ws2.Range("D3").Value = 0.0684408346574217
$ws2.Range("E3").Value = 0.169494028570388
$ws2.Range("F3").Value = 948

# Row 4 (completefam / income)
$ws2.Range("C4").Value = 0.151887435239287
$ws2.Range("D4").Value = 0.0900915315987757
$ws2.Range("E4").Value = 0.213683338879799
$ws2.Range("F4").Value = 948
